# Add a new "T4: 17/3/2020" data column (F) to the COVID history sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column.
$ws.Range("F1").Value = "T4: 17/3/2020"

# Per-department counts for the new date column (rows 2-19).
$counts = @(1, 1, 0, 0, 0, 1, 0, 6, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $counts.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $counts[$i]
}

# Totals row: Sum of the new column.
$ws.Range("F20").Formula = "=SUM(F2:F19)"

# "El Paraiso" (row 8) count is highlighted with an underlined font.
$xlUnderlineStyleSingle = 2
$ws.Range("F8").Font.Underline = $xlUnderlineStyleSingle

# Leave the active selection on the newly edited cell, matching the
# workbook state captured by the author's edit.
$ws.Range("F8").Select()
